# Build 아이템 총 정리
# Update the SpritePath / PrefabPath columns (H, I) for each build item row
# so that each item (작업대/Table, 오두막/Cabin, 침대/Bed) points at its own
# dedicated sprite/prefab asset instead of all sharing the old "well" paths.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BuildItemData")

# Row 2 - 작업대 (workbench/table)
$ws.Range("H2").Value = "Sprites/Table"
$ws.Range("I2").Value = "Prefabs/BuildItemPrefabs/Table"

# Row 3 - 오두막 (cabin)
$ws.Range("H3").Value = "Sprites/Cabin"
$ws.Range("I3").Value = "Prefabs/BuildItemPrefabs/Cabin"

# Row 4 - 침대 (bed)
$ws.Range("H4").Value = "Sprites/Bed"
$ws.Range("I4").Value = "Prefabs/BuildItemPrefabs/Bed"

# Reflect the last active cell selection as seen in the saved file
$ws.Range("E6").Select()

$wb.Save()
